$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================
# Step 1: Swap rows 118 and 119 (columns B:AC), keep A as-is
# =========================================================
$ws.Range("B118").Value = 7013409
$ws.Range("C118").Value = "Uruguay Primera División"
$ws.Range("D118").Value = "Uruguay Clausura"
$ws.Range("E118").Value = 45267.70833333334
$ws.Range("F118").Value = "Nacional De Football"
$ws.Range("G118").Value = "Torque"
$ws.Range("H118").Value = 1
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = "D"
$ws.Range("K118").Value = 1.666
$ws.Range("L118").Value = 3.9
$ws.Range("M118").Value = 4.5
$ws.Range("N118").Value = 1.615
$ws.Range("O118").Value = 4
$ws.Range("P118").Value = 4.75
$ws.Range("Q118").Value = -0.75
$ws.Range("R118").Value = 1.8
$ws.Range("S118").Value = 2.05
$ws.Range("T118").Value = 2.75
$ws.Range("U118").Value = 1.95
$ws.Range("V118").Value = 1.9
$ws.Range("W118").Value = -1
$ws.Range("X118").Value = 3
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = 1.05
$ws.Range("AB118").Value = -1
$ws.Range("AC118").Value = 0.8999999999999999

$ws.Range("B119").Value = 7013702
$ws.Range("C119").Value = "Uruguay Primera División"
$ws.Range("D119").Value = "Uruguay Clausura"
$ws.Range("E119").Value = 45267.70833333334
$ws.Range("F119").Value = "Defensor Sporting"
$ws.Range("G119").Value = "Danubio"
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 2
$ws.Range("J119").Value = "A"
$ws.Range("K119").Value = 1.8
$ws.Range("L119").Value = 3.6
$ws.Range("M119").Value = 4.2
$ws.Range("N119").Value = 1.8
$ws.Range("O119").Value = 3.6
$ws.Range("P119").Value = 4.2
$ws.Range("Q119").Value = -0.75
$ws.Range("R119").Value = 2.05
$ws.Range("S119").Value = 1.8
$ws.Range("T119").Value = 2.25
$ws.Range("U119").Value = 1.85
$ws.Range("V119").Value = 2
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = 3.2
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = 0.8
$ws.Range("AB119").Value = 0.8999999999999999
$ws.Range("AC119").Value = 0.5

# =========================================================
# Step 2: Apply formatting (style) for new rows 185-190
#   column A -> bold/border style (same as A183)
#   column E -> date style (same as E183)
# Do this BEFORE writing values, sourced from an unaffected row
# =========================================================
$ws.Range("A183").Copy()
$ws.Range("A185").PasteSpecial(-4122)
$ws.Range("A186").PasteSpecial(-4122)
$ws.Range("A187").PasteSpecial(-4122)
$ws.Range("A188").PasteSpecial(-4122)
$ws.Range("A189").PasteSpecial(-4122)
$ws.Range("A190").PasteSpecial(-4122)
$ws.Range("E183").Copy()
$ws.Range("E185").PasteSpecial(-4122)
$ws.Range("E186").PasteSpecial(-4122)
$ws.Range("E187").PasteSpecial(-4122)
$ws.Range("E188").PasteSpecial(-4122)
$ws.Range("E189").PasteSpecial(-4122)
$ws.Range("E190").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# =========================================================
# Step 3: Move old row 184 contents down into row 185
#   (old match 8081163 shifts from id 182/row184 to id 183/row185)
# =========================================================
$ws.Range("A185").Value = 183
$ws.Range("B185").Value = 8081163
$ws.Range("C185").Value = "Uruguay Primera División"
$ws.Range("D185").Value = "Uruguay Apertura"
$ws.Range("E185").Value = 45401.5625
$ws.Range("F185").Value = "Racing Club de Montevideo"
$ws.Range("G185").Value = "Cerro"
$ws.Range("K185").Value = 2.3
$ws.Range("L185").Value = 3.2
$ws.Range("M185").Value = 3.2
$ws.Range("N185").Value = 2.1
$ws.Range("O185").Value = 3.3
$ws.Range("P185").Value = 3.6
$ws.Range("Q185").Value = -0.25
$ws.Range("R185").Value = 1.8
$ws.Range("S185").Value = 2.05
$ws.Range("T185").Value = 2.25
$ws.Range("U185").Value = 1.9
$ws.Range("V185").Value = 1.95
$ws.Range("W185").Value = 0
$ws.Range("X185").Value = 0
$ws.Range("Y185").Value = 0
$ws.Range("Z185").Value = 0
$ws.Range("AA185").Value = 0

# =========================================================
# Step 4: Overwrite row 184 with the new match data
# =========================================================
$ws.Range("A184").Value = 182
$ws.Range("B184").Value = 8050913
$ws.Range("C184").Value = "Uruguay Primera División"
$ws.Range("D184").Value = "Uruguay Apertura"
$ws.Range("E184").Value = 45398.85416666666
$ws.Range("F184").Value = "Miramar Misiones"
$ws.Range("G184").Value = "Nacional De Football"
$ws.Range("H184").Value = 1
$ws.Range("I184").Value = 2
$ws.Range("J184").Value = "A"
$ws.Range("K184").Value = 5
$ws.Range("L184").Value = 3.5
$ws.Range("M184").Value = 1.727
$ws.Range("N184").Value = 5.75
$ws.Range("O184").Value = 3.75
$ws.Range("P184").Value = 1.6
$ws.Range("Q184").Value = 1
$ws.Range("R184").Value = 1.75
$ws.Range("S184").Value = 2.05
$ws.Range("T184").Value = 2.25
$ws.Range("U184").Value = 1.95
$ws.Range("V184").Value = 1.9
$ws.Range("W184").Value = -1
$ws.Range("X184").Value = -1
$ws.Range("Y184").Value = 0.6000000000000001
$ws.Range("Z184").Value = 0
$ws.Range("AA184").Value = 0
$ws.Range("AB184").Value = 0.95
$ws.Range("AC184").Value = -1

# =========================================================
# Step 5: Add new rows 186-190
# =========================================================
# row 186
$ws.Range("A186").Value = 184
$ws.Range("B186").Value = 8081435
$ws.Range("C186").Value = "Uruguay Primera División"
$ws.Range("D186").Value = "Uruguay Apertura"
$ws.Range("E186").Value = 45401.64583333334
$ws.Range("F186").Value = "Liverpool Montevideo"
$ws.Range("G186").Value = "Defensor Sporting"
$ws.Range("K186").Value = 2.375
$ws.Range("L186").Value = 3.3
$ws.Range("M186").Value = 3
$ws.Range("N186").Value = 2.625
$ws.Range("O186").Value = 3.3
$ws.Range("P186").Value = 2.7
$ws.Range("Q186").Value = 0
$ws.Range("R186").Value = 1.925
$ws.Range("S186").Value = 1.925
$ws.Range("T186").Value = 2.25
$ws.Range("U186").Value = 1.925
$ws.Range("V186").Value = 1.925
$ws.Range("W186").Value = 0
$ws.Range("X186").Value = 0
$ws.Range("Y186").Value = 0
$ws.Range("Z186").Value = 0
$ws.Range("AA186").Value = 0

# row 187
$ws.Range("A187").Value = 185
$ws.Range("B187").Value = 8081162
$ws.Range("C187").Value = "Uruguay Primera División"
$ws.Range("D187").Value = "Uruguay Apertura"
$ws.Range("E187").Value = 45402.41666666666
$ws.Range("F187").Value = "Danubio"
$ws.Range("G187").Value = "Cerro Largo"
$ws.Range("K187").Value = 2.3
$ws.Range("L187").Value = 3
$ws.Range("M187").Value = 3.4
$ws.Range("N187").Value = 2.15
$ws.Range("O187").Value = 3
$ws.Range("P187").Value = 3.75
$ws.Range("Q187").Value = -0.25
$ws.Range("R187").Value = 1.825
$ws.Range("S187").Value = 2.025
$ws.Range("T187").Value = 2
$ws.Range("U187").Value = 1.925
$ws.Range("V187").Value = 1.925
$ws.Range("W187").Value = 0
$ws.Range("X187").Value = 0
$ws.Range("Y187").Value = 0
$ws.Range("Z187").Value = 0
$ws.Range("AA187").Value = 0

# row 188
$ws.Range("A188").Value = 186
$ws.Range("B188").Value = 8081144
$ws.Range("C188").Value = "Uruguay Primera División"
$ws.Range("D188").Value = "Uruguay Apertura"
$ws.Range("E188").Value = 45402.625
$ws.Range("F188").Value = "Boston River"
$ws.Range("G188").Value = "Penarol"
$ws.Range("K188").Value = 4.75
$ws.Range("L188").Value = 3.75
$ws.Range("M188").Value = 1.727
$ws.Range("N188").Value = 3.75
$ws.Range("O188").Value = 3.5
$ws.Range("P188").Value = 2
$ws.Range("Q188").Value = 0.5
$ws.Range("R188").Value = 1.85
$ws.Range("S188").Value = 2
$ws.Range("T188").Value = 2.25
$ws.Range("U188").Value = 1.975
$ws.Range("V188").Value = 1.875
$ws.Range("W188").Value = 0
$ws.Range("X188").Value = 0
$ws.Range("Y188").Value = 0
$ws.Range("Z188").Value = 0
$ws.Range("AA188").Value = 0

# row 189
$ws.Range("A189").Value = 187
$ws.Range("B189").Value = 8081249
$ws.Range("C189").Value = "Uruguay Primera División"
$ws.Range("D189").Value = "Uruguay Apertura"
$ws.Range("E189").Value = 45402.75
$ws.Range("F189").Value = "Nacional De Football"
$ws.Range("G189").Value = "Rampla Juniors"
$ws.Range("K189").Value = 1.444
$ws.Range("L189").Value = 4
$ws.Range("M189").Value = 8.5
$ws.Range("N189").Value = 1.3
$ws.Range("O189").Value = 4.5
$ws.Range("P189").Value = 12
$ws.Range("Q189").Value = -1.5
$ws.Range("R189").Value = 1.975
$ws.Range("S189").Value = 1.875
$ws.Range("T189").Value = 2.5
$ws.Range("U189").Value = 1.975
$ws.Range("V189").Value = 1.875
$ws.Range("W189").Value = 0
$ws.Range("X189").Value = 0
$ws.Range("Y189").Value = 0
$ws.Range("Z189").Value = 0
$ws.Range("AA189").Value = 0

# row 190
$ws.Range("A190").Value = 188
$ws.Range("B190").Value = 8081250
$ws.Range("C190").Value = "Uruguay Primera División"
$ws.Range("D190").Value = "Uruguay Apertura"
$ws.Range("E190").Value = 45402.85416666666
$ws.Range("F190").Value = "Deportivo Maldonado"
$ws.Range("G190").Value = "Miramar Misiones"
$ws.Range("K190").Value = 2.2
$ws.Range("L190").Value = 3.3
$ws.Range("M190").Value = 3.3
$ws.Range("N190").Value = 2.1
$ws.Range("O190").Value = 3.4
$ws.Range("P190").Value = 3.5
$ws.Range("Q190").Value = -0.25
$ws.Range("R190").Value = 1.825
$ws.Range("S190").Value = 2.025
$ws.Range("T190").Value = 2.25
$ws.Range("U190").Value = 1.975
$ws.Range("V190").Value = 1.875
$ws.Range("W190").Value = 0
$ws.Range("X190").Value = 0
$ws.Range("Y190").Value = 0
$ws.Range("Z190").Value = 0
$ws.Range("AA190").Value = 0
